$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '256.30'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-0.42%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '26.85'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-1.00%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '4.537'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-4.48%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05880'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.609'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-0.89%'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8501'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-2.29%'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-1.79%'
$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = 'One'
$ws.Range('C9').NumberFormat = '@'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0006060'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-0.19%'
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1379'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-2.02%'
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.04487'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '23.90%'
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07025'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-1.98%'
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03067'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-3.26%'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09102'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-1.50%'
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001527'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-1.57%'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.006139'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '1.77%'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.480'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-0.05%'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.165'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-0.45%'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-1.65%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1285'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-1.61%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.911'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '2.56%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04273'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '1.42%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001222'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '0.03%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004300'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-4.43%'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-0.02%'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '2.02%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03802'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-0.25%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006250'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '56.93%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1098'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-0.29%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002200'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-4.36%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.01383'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '30.95%'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-2.75%'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-0.02%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05379'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-39.25%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.2526'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '10,678.75%'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.02%'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.02%'
